$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format numeric-looking price cells as Text so Excel keeps them
# as literal strings (matching the source feed's text formatting)
# instead of auto-converting to the Number type.
$textCells = @("D5", "D6", "D10", "D12", "D13", "D17", "D19", "D20", "D21", "D22", "D25", "D27", "D29", "D32", "D35", "D36", "D37", "D38", "D40", "D42", "D43", "D46", "D47", "D48", "D49")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "72.357.75"
$ws.Range("E2").Value = "  -0.21%  "
$ws.Range("D3").Value = "2.642.04"
$ws.Range("E3").Value = "  -1.38%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "583.33"
$ws.Range("E5").Value = "  -3.28%  "
$ws.Range("D6").Value = "174.93"
$ws.Range("E6").Value = "  -2.08%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  -1.13%  "
$ws.Range("D9").Value = "2.642.01"
$ws.Range("E9").Value = "  -1.38%  "
$ws.Range("D10").Value = "0.171"
$ws.Range("E10").Value = "  -0.31%  "
$ws.Range("E11").Value = "  +0.92%  "
$ws.Range("D12").Value = "0.356"
$ws.Range("E12").Value = "  -0.24%  "
$ws.Range("D13").Value = "4.92"
$ws.Range("E13").Value = "  -2.68%  "
$ws.Range("D14").Value = "3.123.91"
$ws.Range("E14").Value = "  -1.37%  "
$ws.Range("E15").Value = "  -0.77%  "
$ws.Range("D16").Value = "72.207.42"
$ws.Range("E16").Value = "  -0.24%  "
$ws.Range("D17").Value = "25.86"
$ws.Range("E17").Value = "  -1.93%  "
$ws.Range("D18").Value = "2.646.65"
$ws.Range("E18").Value = "  -1.08%  "
$ws.Range("D19").Value = "8.42"
$ws.Range("E19").Value = "  +5.00%  "
$ws.Range("D20").Value = "12.11"
$ws.Range("E20").Value = "  +1.57%  "
$ws.Range("D21").Value = "374.25"
$ws.Range("E21").Value = "  +0.25%  "
$ws.Range("D22").Value = "4.12"
$ws.Range("E22").Value = "  -1.46%  "
$ws.Range("E23").Value = "  -0.28%  "
$ws.Range("E24").Value = "  +0.07%  "
$ws.Range("D25").Value = "70.78"
$ws.Range("E25").Value = "  -2.33%  "
$ws.Range("E26").Value = "  -2.36%  "
$ws.Range("D27").Value = "9.52"
$ws.Range("E27").Value = "  -3.76%  "
$ws.Range("D28").Value = "2.778.23"
$ws.Range("E28").Value = "  -1.31%  "
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").Value = "  -0.11%  "
$ws.Range("D30").Value = "0.0₃0950"
$ws.Range("E30").Value = "  +0.40%  "
$ws.Range("E31").Value = "  -1.66%  "
$ws.Range("D32").Value = "494.93"
$ws.Range("E32").Value = "  -4.61%  "
$ws.Range("E33").Value = "  -2.76%  "
$ws.Range("E34").Value = "  -1.89%  "
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  -0.04%  "
$ws.Range("D36").Value = "162.68"
$ws.Range("E36").Value = "  -1.18%  "
$ws.Range("B37").Value = "Kaspa"
$ws.Range("C37").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D37").Value = "0.115"
$ws.Range("E37").Value = "  +4.71%  "
$ws.Range("B38").Value = "EthereumClassic"
$ws.Range("C38").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D38").Value = "19.17"
$ws.Range("E38").Value = "  -1.89%  "
$ws.Range("E39").Value = "  -1.55%  "
$ws.Range("D40").Value = "1.36"
$ws.Range("E40").Value = "  -2.12%  "
$ws.Range("E41").Value = "  -0.03%  "
$ws.Range("B42").Value = "dogwifhat"
$ws.Range("C42").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D42").Value = "2.57"
$ws.Range("E42").Value = "  -1.11%  "
$ws.Range("B43").Value = "Stacks"
$ws.Range("C43").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D43").Value = "1.72"
$ws.Range("E43").Value = "  -6.57%  "
$ws.Range("E44").Value = "  -3.23%  "
$ws.Range("E45").Value = "  -2.29%  "
$ws.Range("D46").Value = "39.00"
$ws.Range("E46").Value = "  -0.59%  "
$ws.Range("D47").Value = "152.03"
$ws.Range("E47").Value = "  -1.29%  "
$ws.Range("D48").Value = "3.65"
$ws.Range("E48").Value = "  -2.54%  "
$ws.Range("D49").Value = "0.545"
$ws.Range("E49").Value = "  -0.69%  "
$ws.Range("E50").Value = "  -2.97%  "
$ws.Range("E51").Value = "  -0.93%  "
